# Update "Metadata" sheet: refresh the "Last Updated" timestamp (A2).
$wb = $excel.ActiveWorkbook
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("A2").Value = "30 Oct 2025, 10:25 AM"

# Update "distance from Dma50" sheet: refresh Distance-From-Sma50 values
# (column C) for rows 2-30, and apply the row 10/11, 18/19 and 26/27
# stock-name (column B) swaps that came along with the refreshed ranking.
$wsDma = $wb.Worksheets.Item("distance from Dma50")

$wsDma.Range("C2").Value = 9.7285
$wsDma.Range("C3").Value = 7.6214
$wsDma.Range("C5").Value = 5.2911
$wsDma.Range("C6").Value = 5.141
$wsDma.Range("C7").Value = 4.8883
$wsDma.Range("C8").Value = 4.4614
$wsDma.Range("C9").Value = 4.3319

$wsDma.Range("B10").Value = "NIFTYMIDCAP50"
$wsDma.Range("C10").Value = 3.5794
$wsDma.Range("B11").Value = "NIFTYFINSERVICE"
$wsDma.Range("C11").Value = 3.5298

$wsDma.Range("C12").Value = 3.5248
$wsDma.Range("C13").Value = 3.138
$wsDma.Range("C14").Value = 3.1316
$wsDma.Range("C15").Value = 3.0367
$wsDma.Range("C16").Value = 3.0091
$wsDma.Range("C17").Value = 2.8002

$wsDma.Range("B18").Value = "NIFTYCPSE"
$wsDma.Range("C18").Value = 2.5478
$wsDma.Range("B19").Value = "CNXSMALLCAP"
$wsDma.Range("C19").Value = 2.5465

$wsDma.Range("C20").Value = 2.3384
$wsDma.Range("C21").Value = 2.2157
$wsDma.Range("C22").Value = 1.3069
$wsDma.Range("C23").Value = 1.3047
$wsDma.Range("C24").Value = 1.0328
$wsDma.Range("C25").Value = 1.0259

$wsDma.Range("B26").Value = "NIFTYGROWSECT15"
$wsDma.Range("C26").Value = 0.8645
$wsDma.Range("B27").Value = "NIFTYFMCG"
$wsDma.Range("C27").Value = 0.8423

$wsDma.Range("C28").Value = 0.4189
$wsDma.Range("C29").Value = -0.1776
$wsDma.Range("C30").Value = -2.099
